# Auto-generated Excel COM-interop script
# Applies numeric corrections to the Leve profit-tracking columns (H-N)
# across multiple worksheets, per the scheduled runner update.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 41666780
$ws.Range("I9").Value = 62500100
$ws.Range("K9").Value = 62500100
$ws.Range("M9").Value = -62499931

$ws.Range("H58").Value = 253.6
$ws.Range("J58").Value = 388.66666
$ws.Range("L58").Value = 1165.99998
$ws.Range("N58").Value = -1465.99998

$ws.Range("H63").Value = 54999
$ws.Range("I63").Value = 54999
$ws.Range("K63").Value = 54999
$ws.Range("M63").Value = -54375

$ws.Range("H66").Value = 54999
$ws.Range("I66").Value = 54999
$ws.Range("K66").Value = 164997
$ws.Range("M66").Value = -161877

$ws.Range("H86").Value = 10054111
$ws.Range("I86").Value = 4468.4
$ws.Range("K86").Value = 4468.4
$ws.Range("M86").Value = -3345.4

$ws.Range("H88").Value = 2455.4443
$ws.Range("I88").Value = 800
$ws.Range("J88").Value = 2662.375
$ws.Range("K88").Value = 800
$ws.Range("L88").Value = 2662.375
$ws.Range("M88").Value = -394
$ws.Range("N88").Value = -3474.375

$ws.Range("H89").Value = 10054111
$ws.Range("I89").Value = 4468.4
$ws.Range("K89").Value = 22342
$ws.Range("M89").Value = -16726

$ws.Range("H91").Value = 2455.4443
$ws.Range("I91").Value = 800
$ws.Range("J91").Value = 2662.375
$ws.Range("K91").Value = 800
$ws.Range("L91").Value = 2662.375
$ws.Range("M91").Value = 604
$ws.Range("N91").Value = -5470.375

$ws.Range("H113").Value = 200002610
$ws.Range("I113").Value = 250002260
$ws.Range("J113").Value = 4000
$ws.Range("K113").Value = 250002260
$ws.Range("L113").Value = 4000
$ws.Range("M113").Value = -249999006
$ws.Range("N113").Value = -10508

$ws.Range("H137").Value = 3928.8
$ws.Range("I137").Value = 3028.2
$ws.Range("K137").Value = 9084.599999999999
$ws.Range("M137").Value = -6534.599999999999

$ws.Range("H138").Value = 4670.5
$ws.Range("I138").Value = 4056.75
$ws.Range("J138").Value = 4962.7617
$ws.Range("K138").Value = 12170.25
$ws.Range("L138").Value = 14888.2851
$ws.Range("M138").Value = -7030.25
$ws.Range("N138").Value = -25168.2851

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2000
$ws.Range("I2").Value = 2000
$ws.Range("J2").Value = 2000
$ws.Range("K2").Value = 2000
$ws.Range("L2").Value = 2000
$ws.Range("M2").Value = -1887
$ws.Range("N2").Value = -2226

$ws.Range("H32").Value = 16374.284
$ws.Range("I32").Value = 11934.07
$ws.Range("J32").Value = 47899.8
$ws.Range("K32").Value = 11934.07
$ws.Range("L32").Value = 47899.8
$ws.Range("M32").Value = -11647.07
$ws.Range("N32").Value = -48473.8

$ws.Range("H116").Value = 2000
$ws.Range("I116").Value = 2000
$ws.Range("J116").Value = 2000
$ws.Range("K116").Value = 2000
$ws.Range("L116").Value = 2000
$ws.Range("M116").Value = 294
$ws.Range("N116").Value = -6588

$ws.Range("H122").Value = 5242.1816
$ws.Range("I122").Value = 2333
$ws.Range("K122").Value = 6999
$ws.Range("M122").Value = -4549

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2000
$ws.Range("I3").Value = 2000
$ws.Range("J3").Value = 2000
$ws.Range("K3").Value = 2000
$ws.Range("L3").Value = 2000
$ws.Range("M3").Value = -1886
$ws.Range("N3").Value = -2228

$ws.Range("H69").Value = 34999.5
$ws.Range("I69").Value = 34999.5
$ws.Range("K69").Value = 34999.5
$ws.Range("M69").Value = -34188.5

$ws.Range("H72").Value = 34999.5
$ws.Range("I72").Value = 34999.5
$ws.Range("K72").Value = 104998.5
$ws.Range("M72").Value = -100942.5

$ws.Range("H99").Value = 3633
$ws.Range("I99").Value = 950
$ws.Range("K99").Value = 950
$ws.Range("M99").Value = 548

$ws.Range("H134").Value = 3309.4722
$ws.Range("I134").Value = 2927.2856
$ws.Range("J134").Value = 4647.125
$ws.Range("K134").Value = 8781.856800000001
$ws.Range("L134").Value = 13941.375
$ws.Range("M134").Value = -6246.856800000001
$ws.Range("N134").Value = -19011.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3966.3215
$ws.Range("I31").Value = 3005
$ws.Range("J31").Value = 5452
$ws.Range("K31").Value = 3005
$ws.Range("L31").Value = 5452
$ws.Range("M31").Value = -2710
$ws.Range("N31").Value = -6042

$ws.Range("H34").Value = 3966.3215
$ws.Range("I34").Value = 3005
$ws.Range("J34").Value = 5452
$ws.Range("K34").Value = 3005
$ws.Range("L34").Value = 5452
$ws.Range("M34").Value = -2803
$ws.Range("N34").Value = -5856

$ws.Range("H122").Value = 5255.5
$ws.Range("I122").Value = 513
$ws.Range("K122").Value = 1539
$ws.Range("M122").Value = 911

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 44285.715
$ws.Range("J37").Value = 44285.715
$ws.Range("L37").Value = 132857.145
$ws.Range("N37").Value = -133081.145

$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()

$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()

$ws.Range("H80").Value = 6057.9565
$ws.Range("I80").Value = 5843.1177
$ws.Range("K80").Value = 17529.3531
$ws.Range("M80").Value = -16593.3531

$ws.Range("H83").Value = 6057.9565
$ws.Range("I83").Value = 5843.1177
$ws.Range("K83").Value = 52588.0593
$ws.Range("M83").Value = -47908.0593

$ws.Range("H121").Value = 13441.25
$ws.Range("I121").Value = 14782.857
$ws.Range("K121").Value = 44348.571
$ws.Range("M121").Value = -43038.571

$ws.Range("H137").Value = 4777.4
$ws.Range("J137").Value = 5380.9165
$ws.Range("L137").Value = 16142.7495
$ws.Range("N137").Value = -26342.7495

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 84997.5
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 84997.5
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 84997.5
$ws.Range("N63").Value = -86369.5
$ws.Range("M63").ClearContents()

$ws.Range("H66").Value = 84997.5
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 84997.5
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 254992.5
$ws.Range("N66").Value = -261856.5
$ws.Range("M66").ClearContents()

$ws.Range("H80").Value = 2151.2942
$ws.Range("J80").Value = 2698.4443
$ws.Range("L80").Value = 2698.4443
$ws.Range("N80").Value = -4694.4443

$ws.Range("H83").Value = 2151.2942
$ws.Range("J83").Value = 2698.4443
$ws.Range("L83").Value = 13492.2215
$ws.Range("N83").Value = -23476.2215

$ws.Range("H102").Value = 2209.1956
$ws.Range("I102").Value = 959.1667
$ws.Range("J102").Value = 3572.8635
$ws.Range("K102").Value = 959.1667
$ws.Range("L102").Value = 3572.8635
$ws.Range("M102").Value = 662.8333
$ws.Range("N102").Value = -6816.863499999999

$ws.Range("H122").Value = 399968
$ws.Range("I122").Value = 796222.1
$ws.Range("K122").Value = 2388666.3
$ws.Range("M122").Value = -2386216.3

$ws.Range("H126").Value = 9492.519
$ws.Range("I126").Value = 19857
$ws.Range("J126").Value = 5864.95
$ws.Range("K126").Value = 59571
$ws.Range("L126").Value = 17594.85
$ws.Range("M126").Value = -57101
$ws.Range("N126").Value = -22534.85

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 550
$ws.Range("I22").Value = 500
$ws.Range("J22").Value = 600
$ws.Range("K22").Value = 500
$ws.Range("L22").Value = 600
$ws.Range("M22").Value = -205
$ws.Range("N22").Value = -1190

$ws.Range("H27").Value = 550
$ws.Range("I27").Value = 500
$ws.Range("J27").Value = 600
$ws.Range("K27").Value = 500
$ws.Range("L27").Value = 600
$ws.Range("M27").Value = -393
$ws.Range("N27").Value = -814

$ws.Range("H40").Value = 4576.6313
$ws.Range("I40").Value = 3494.875
$ws.Range("J40").Value = 5363.364
$ws.Range("K40").Value = 3494.875
$ws.Range("L40").Value = 5363.364
$ws.Range("M40").Value = -3358.875
$ws.Range("N40").Value = -5635.364

$ws.Range("H46").Value = 266528.34
$ws.Range("I46").Value = 2859
$ws.Range("J46").Value = 360695.97
$ws.Range("K46").Value = 2859
$ws.Range("L46").Value = 360695.97
$ws.Range("M46").Value = -2671
$ws.Range("N46").Value = -361071.97

$ws.Range("H68").Value = 6446.615
$ws.Range("I68").Value = 6446.615
$ws.Range("K68").Value = 6446.615
$ws.Range("M68").Value = -5697.615

$ws.Range("H71").Value = 6446.615
$ws.Range("I71").Value = 6446.615
$ws.Range("K71").Value = 32233.075
$ws.Range("M71").Value = -28489.075

$ws.Range("H100").Value = 3180
$ws.Range("I100").Value = 1966.6666
$ws.Range("J100").Value = 5000
$ws.Range("K100").Value = 1966.6666
$ws.Range("L100").Value = 5000
$ws.Range("M100").Value = -1425.6666
$ws.Range("N100").Value = -6082

$ws.Range("H132").Value = 131463.86
$ws.Range("J132").Value = 6251.231
$ws.Range("L132").Value = 18753.693
$ws.Range("N132").Value = -23813.693

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 5986.0347
$ws.Range("I113").Value = 7223.722
$ws.Range("K113").Value = 21671.166
$ws.Range("M113").Value = -19501.166

$ws.Range("H133").Value = 116883.85
$ws.Range("J133").Value = 120930.37
$ws.Range("L133").Value = 120930.37
$ws.Range("N133").Value = -131050.37

$ws.Range("H136").Value = 6883.067
$ws.Range("I136").Value = 6453.684
$ws.Range("K136").Value = 19361.052
$ws.Range("M136").Value = -16811.052
